# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 14 (pushing the existing
# rows 14:60 down to 15:61), and the new row is populated with the
# corresponding market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 14:60 down by one row, opening up a blank row 14.
$ws.Rows(14).Insert()

# Fill the newly inserted row 14 with the new record's data.
$fecha = Get-Date -Year 2022 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(14, 1).Value = 7
$ws.Cells.Item(14, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14, 3).Value = "Ñuble"
$ws.Cells.Item(14, 4).Value = $fecha
$ws.Cells.Item(14, 5).Value = 16
$ws.Cells.Item(14, 6).Value = 100112031
$ws.Cells.Item(14, 7).Value = "Poroto verde"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 120
$ws.Cells.Item(14, 11).Value = 25000
$ws.Cells.Item(14, 12).Value = 26000
$ws.Cells.Item(14, 13).Value = 25500
$ws.Cells.Item(14, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región del Maule"
$ws.Cells.Item(14, 16).Value = 1020
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"
